# "Revert 'Revert CCS files to same as master branch'" -------------------------
# Restores the CCS Percentages-by-Entity data to the pre-revert (master)
# state on the two sheets that carry the CCS entity-share tables:
#   - CPbE-FoCSbS      (CCS Percentage by Entity - Fraction of CO2 Sequestration by Sector)
#   - CPbE-FoESCbES    (CCS Percentage by Entity - Fraction of Electricity Sector CCS by Energy Source)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "CPbE-FoCSbS": electricity-sector / industry-sector split flips from
# (0%, 100%) back to (100%, 0%), and the formulas driving it are flattened to
# hard values (row 2), while row 3 keeps its "=$B$3" style formulas but now
# recomputes to 0 because B3 itself becomes a literal 0.
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("CPbE-FoCSbS")

# Row 2 ("electricity sector"): was "=About!$I$24" -> 0, now a literal 1.
$wsA.Range("B2:AM2").Value = 1

# Row 3 ("industry sector"): B3 was "=About!$I$23" -> 1, now a literal 0.
# C3:AM3 keep referencing $B$3 (directly or via the shared formula), so they
# fall back to 0 automatically once B3 is rewritten.
$wsA.Range("B3").Value = 0

$wsA.Activate()
$wsA.Range("B4").Select()

# ---------------------------------------------------------------------------
# Sheet "CPbE-FoESCbES": every energy-source row collapses to hard-coded
# values - zero for every source except "natural gas peaker" (row 12), which
# takes the full 100% share. Row 13 ("lignite") keeps formulas, but each
# cell now just points at its neighbor to the right (a cascading, unseeded
# chain), which evaluates to 0 all the way across.
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("CPbE-FoESCbES")

$wsB.Range("B2:AM11").Value = 0
$wsB.Range("B12:AM12").Value = 1
$wsB.Range("B14:AM14").Value = 0

$cols = @("B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z", `
          "AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ","AK","AL","AM")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $col = $cols[$i]
    if ($i -lt $cols.Length - 1) {
        $nextCol = $cols[$i + 1]
    } else {
        $nextCol = "AN"
    }
    $wsB.Range($col + "13").Formula = "=" + $nextCol + "13"
}

$wsB.Activate()
$wsB.Range("B14:AM14").Select()
